# Updates Kraken Profits leve-profit figures (currentAveragePrice* / LevePrice* / LeveProfit*
# columns H:N) for the rows touched by the scheduled price refresh, one worksheet per job class.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether
$ws.Range("H15").Value = 1612.1538
$ws.Range("I15").Value = 1612.1538
$ws.Range("K15").Value = 4836.4614
$ws.Range("M15").Value = -4667.4614

# Row 51: A Bile Business
$ws.Range("H51").Value = 6666.6665
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 7000
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 7000
$ws.Range("M51").Value = -4516
$ws.Range("N51").Value = -7968

# Row 53: No Accounting for Waste
$ws.Range("H53").Value = 810.0909
$ws.Range("I53").Value = 92.625
$ws.Range("J53").Value = 2723.3333
$ws.Range("K53").Value = 92.625
$ws.Range("L53").Value = 2723.3333
$ws.Range("M53").Value = 544.375
$ws.Range("N53").Value = -3997.3333

# Row 88: The Grave of Hemlock Groves
$ws.Range("H88").Value = 4104.6
$ws.Range("I88").Value = 3533.3333
$ws.Range("J88").Value = 4961.5
$ws.Range("K88").Value = 3533.3333
$ws.Range("L88").Value = 4961.5
$ws.Range("M88").Value = -3127.3333
$ws.Range("N88").Value = -5773.5

# Row 91: Dappling the Highlands (L)
$ws.Range("H91").Value = 4104.6
$ws.Range("I91").Value = 3533.3333
$ws.Range("J91").Value = 4961.5
$ws.Range("K91").Value = 3533.3333
$ws.Range("L91").Value = 4961.5
$ws.Range("M91").Value = -2129.3333
$ws.Range("N91").Value = -7769.5

# Row 94: Magic Beans
$ws.Range("H94").Value = 1494
$ws.Range("I94").Value = 1494
$ws.Range("K94").Value = 1494
$ws.Range("M94").Value = -1043

# Row 113: Amaro Kart
$ws.Range("H113").Value = 966.6667
$ws.Range("I113").Value = 950
$ws.Range("K113").Value = 950
$ws.Range("M113").Value = 2304

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 673
$ws.Range("I2").Value = 618.5
$ws.Range("K2").Value = 618.5
$ws.Range("M2").Value = -505.5

# Row 116: No Scope
$ws.Range("H116").Value = 673
$ws.Range("I116").Value = 618.5
$ws.Range("K116").Value = 618.5
$ws.Range("M116").Value = 1675.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 673
$ws.Range("I3").Value = 618.5
$ws.Range("K3").Value = 618.5
$ws.Range("M3").Value = -504.5

# Row 132: Always Be Prepaired
$ws.Range("H132").Value = 99780
$ws.Range("J132").Value = 99780
$ws.Range("L132").Value = 99780
$ws.Range("N132").Value = -109900

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value = 2377.75
$ws.Range("I16").Value = 2405.5
$ws.Range("J16").Value = 2350
$ws.Range("K16").Value = 2405.5
$ws.Range("L16").Value = 2350
$ws.Range("M16").Value = -2118.5
$ws.Range("N16").Value = -2924

# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 275
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Row 86: Birch, Please
$ws.Range("H86").Value = 4118.5
$ws.Range("I86").Value = 4118.5
$ws.Range("K86").Value = 4118.5
$ws.Range("M86").Value = -2995.5

# Row 89: Built This City on Blocks and Soul (L)
$ws.Range("H89").Value = 4118.5
$ws.Range("I89").Value = 4118.5
$ws.Range("K89").Value = 20592.5
$ws.Range("M89").Value = -14976.5

# Row 113: Patient Patients
$ws.Range("H113").Value = 2377.75
$ws.Range("I113").Value = 2405.5
$ws.Range("J113").Value = 2350
$ws.Range("K113").Value = 2405.5
$ws.Range("L113").Value = 2350
$ws.Range("M113").Value = -235.5
$ws.Range("N113").Value = -6690

$ws = $wb.Worksheets.Item("CUL")
# Row 11: Putting the Squeeze On
$ws.Range("H11").Value = 185.5
$ws.Range("I11").Value = 185.5
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 556.5
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -416.5
$ws.Range("N11").ClearContents()

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# Row 132: More Mezcal
$ws.Range("H132").Value = 2620.6
$ws.Range("I132").Value = 3101
$ws.Range("J132").Value = 699
$ws.Range("K132").Value = 27909
$ws.Range("L132").Value = 6291
$ws.Range("M132").Value = -25379
$ws.Range("N132").Value = -11351

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 7299.706
$ws.Range("I122").Value = 1930.4615
$ws.Range("J122").Value = 24749.75
$ws.Range("K122").Value = 5791.3845
$ws.Range("L122").Value = 74249.25
$ws.Range("M122").Value = -3341.3845
$ws.Range("N122").Value = -79149.25

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad
$ws.Range("H40").Value = 7999.5
$ws.Range("J40").Value = 7999
$ws.Range("L40").Value = 7999
$ws.Range("N40").Value = -8271

# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 5874.625
$ws.Range("I61").Value = 5874.625
$ws.Range("K61").Value = 5874.625
$ws.Range("M61").Value = -5672.625

# Row 113: Peace in Rest
$ws.Range("H113").Value = 5874.625
$ws.Range("I113").Value = 5874.625
$ws.Range("K113").Value = 5874.625
$ws.Range("M113").Value = -3704.625

# Row 122: Hell on Leather
$ws.Range("H122").Value = 4141.1113
$ws.Range("I122").Value = 4474.6
$ws.Range("J122").Value = 3724.25
$ws.Range("K122").Value = 13423.8
$ws.Range("L122").Value = 11172.75
$ws.Range("M122").Value = -10973.8
$ws.Range("N122").Value = -16072.75

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 25290.3
$ws.Range("I132").Value = 26129
$ws.Range("K132").Value = 78387
$ws.Range("M132").Value = -75857

$ws = $wb.Worksheets.Item("WVR")
# Row 2: The Unmentionables
$ws.Range("H2").Value = 2041.625
$ws.Range("I2").Value = 1618.5714
$ws.Range("J2").Value = 5003
$ws.Range("K2").Value = 1618.5714
$ws.Range("L2").Value = 5003
$ws.Range("M2").Value = -1506.5714
$ws.Range("N2").Value = -5227

# Row 4: Not Cool Enough
$ws.Range("H4").Value = 673.9
$ws.Range("I4").Value = 92.5
$ws.Range("J4").Value = 2999.5
$ws.Range("K4").Value = 92.5
$ws.Range("L4").Value = 2999.5
$ws.Range("M4").Value = 20.5
$ws.Range("N4").Value = -3225.5

# Row 33: I'll Be Your Wailer Today
$ws.Range("H33").Value = 22509.5
$ws.Range("I33").Value = 15019
$ws.Range("J33").Value = 30000
$ws.Range("K33").Value = 15019
$ws.Range("L33").Value = 30000
$ws.Range("M33").Value = -14769
$ws.Range("N33").Value = -30500

# Row 36: Put a Lid on It
$ws.Range("H36").Value = 22509.5
$ws.Range("I36").Value = 15019
$ws.Range("J36").Value = 30000
$ws.Range("K36").Value = 15019
$ws.Range("L36").Value = 30000
$ws.Range("M36").Value = -14769
$ws.Range("N36").Value = -30500

# Row 37: Bet You Anything
$ws.Range("H37").Value = 25026
$ws.Range("I37").Value = 25026
$ws.Range("K37").Value = 25026
$ws.Range("M37").Value = -24823

# Row 100: Of Great Import
$ws.Range("H100").Value = 5785.9287
$ws.Range("I100").Value = 6076.923
$ws.Range("J100").Value = 2003
$ws.Range("K100").Value = 12153.846
$ws.Range("L100").Value = 4006
$ws.Range("M100").Value = -11612.846
$ws.Range("N100").Value = -5088

# Row 122: Heavy Armoire
$ws.Range("H122").Value = 202522.2
$ws.Range("I122").Value = 287603.16
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 862809.48
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -860359.48
$ws.Range("N122").Value = -16900

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 1158.6666
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 4814.857
$ws.Range("I132").Value = 4814.857
$ws.Range("K132").Value = 14444.571
$ws.Range("M132").Value = -11914.571
